# Append two new match rows (181 and 182) to Sheet1, mirroring the
# formatting of the last existing data row (180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone formatting from the last existing row into the two new rows ---
$ws.Range("A180:V180").Copy()
$ws.Range("A181:V182").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# "temporada" (column D) is stored as text "2023"; copy the value straight
# from the row above rather than assigning a literal, which Excel would
# otherwise auto-convert to a number.
$ws.Range("D180").Copy()
$ws.Range("D181").PasteSpecial(-4163)
$ws.Range("D180").Copy()
$ws.Range("D182").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Row 181: Aguilas 3 x 1 Envigado ---
$ws.Range("A181").Value = 180
$ws.Range("B181").Value = "colombia"
$ws.Range("C181").Value = "primera-a"
$ws.Range("E181").Value = 45224.04166666666
$ws.Range("F181").Value = "Aguilas"
$ws.Range("G181").Value = 3
$ws.Range("H181").Value = "Envigado"
$ws.Range("I181").Value = 1
$ws.Range("J181").Value = 1.38
$ws.Range("K181").Value = "22/10/2023 22:42"
$ws.Range("L181").Value = 1.48
$ws.Range("M181").Value = "25/10/2023 00:53"
$ws.Range("N181").Value = 4.65
$ws.Range("O181").Value = "22/10/2023 22:42"
$ws.Range("P181").Value = 4.56
$ws.Range("Q181").Value = "25/10/2023 00:53"
$ws.Range("R181").Value = 7.73
$ws.Range("S181").Value = "22/10/2023 22:42"
$ws.Range("T181").Value = 7.07
$ws.Range("U181").Value = "25/10/2023 00:53"
$ws.Range("V181").Value = "https://www.betexplorer.com/football/colombia/primera-a/aguilas-doradas-envigado/GxOE6aes/"

# --- Row 182: Dep. Cali 2 x 0 Jaguares de Cordoba ---
$ws.Range("A182").Value = 181
$ws.Range("B182").Value = "colombia"
$ws.Range("C182").Value = "primera-a"
$ws.Range("E182").Value = 45224.13194444445
$ws.Range("F182").Value = "Dep. Cali"
$ws.Range("G182").Value = 2
$ws.Range("H182").Value = "Jaguares de Cordoba"
$ws.Range("I182").Value = 0
$ws.Range("J182").Value = 1.57
$ws.Range("K182").Value = "22/10/2023 22:42"
$ws.Range("L182").Value = 1.44
$ws.Range("M182").Value = "25/10/2023 03:07"
$ws.Range("N182").Value = 3.78
$ws.Range("O182").Value = "22/10/2023 22:42"
$ws.Range("P182").Value = 4.37
$ws.Range("Q182").Value = "25/10/2023 03:07"
$ws.Range("R182").Value = 6.88
$ws.Range("S182").Value = "22/10/2023 22:42"
$ws.Range("T182").Value = 8.81
$ws.Range("U182").Value = "25/10/2023 03:09"
$ws.Range("V182").Value = "https://www.betexplorer.com/football/colombia/primera-a/dep-cali-jaguares-de-cordoba/lArO87rP/"
